# "Added Delete and View Functionality" - log a new changelog entry on the
# tracking sheet: a date in column A (formatted as yyyy-mm-dd) and a
# description in column B, appended right below the existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 44492
$ws.Range("A2").NumberFormat = "yyyy-mm-dd"
$ws.Range("B2").Value = "Added View Tab Implementation in Personal Organizer"
